$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump the Quantity for row 3 (Capacitor / Ceramic / 22pF) from 43 to 44.
# The Total column (K) is a calculated table column
# (=Table3[[#This Row],[Price]]*Table3[[#This Row],[Quantity]]) and the
# totals row (J21/K21) are SUM() formulas over the table, so both will
# recalc automatically.
$ws.Range("J3").Value = 44

# Move the active selection from H16 to J16, as recorded in the sheetView.
$ws.Range("J16").Select()

$wb.Save()
